$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 = lab06. Mark "link_it" (column C) TRUE and add the new topic (column D).
$ws.Range("C27").Value = $true
$ws.Range("D27").Value = "Webdata in Python"

# Update the active selection to match the author's final cursor position.
$ws.Range("D28").Select()
